$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.230985683306322
$ws.Range("C2").Value = 1.667794583268128
$ws.Range("D2").Value = 3.900430680208489
$ws.Range("E2").Value = 0.496779210170732
$ws.Range("G2").Value = 9.295990156953671
